$d = $word.ActiveDocument

# 1. "...and about 9 dB better sensitivity than the normal 2-minute mode..."
#    -> drop the leading "about " qualifier.
$d.Content.Find.Execute("about 9 dB better sensitivity than ", $true, $false, $false, $false, $false, $true, 1, $false, "9 dB better sensitivity than ", 2) | Out-Null

# 2. Refresh the cached PAGE field result shown in the footer (was stale at "4").
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null
